# Update numeric values in column F ("views"/count column) for the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets, reflecting
# refreshed counters from the upstream data source.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 42
$ws1.Range("F5").Value = 3131
$ws1.Range("F7").Value = 3901
$ws1.Range("F8").Value = 479
$ws1.Range("F9").Value = 979
$ws1.Range("F10").Value = 35

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 42
$ws4.Range("F6").Value = 3131
$ws4.Range("F8").Value = 3901
$ws4.Range("F9").Value = 479
$ws4.Range("F10").Value = 979
$ws4.Range("F11").Value = 35
